$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the reproduction-steps "box" down by one row ------------------
# Row 8 currently closes the bordered box (left+right+bottom border, s=7/17).
# A new step is being added, so that closing border moves down to row 9,
# and row 8 becomes a regular "middle" row (left+right border only, s=5/6),
# matching rows 6-7.
$ws.Range("A8:B8").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)   # xlPasteFormats - push closing style down
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)   # xlPasteFormats - row 8 becomes "middle" style

# --- Title / Description swapped from DF-002 to DF-001 -------------------
$ws.Range("B1").Value = "Se pueden intentar logear DF-001"
$ws.Range("B3").Value = "Al logearse en la cuenta e ingresar mal el password, la pagina identifica que el password es invalido pero te deja intentarlo infinitas veces"

# --- Reproduction steps (rows 6-9) ----------------------------------------
$ws.Range("B7").Value = '3- Click en "Log in"'
$ws.Range("B8").Value = "4-Ingresar mail valido y password random"
$ws.Range("B9").Value = "5-Re intentar infinitas veces."

# --- Cosmetic row/column sizing (best effort) -----------------------------
$ws.Columns.Item(2).ColumnWidth = 75.875
$ws.Rows.Item(24).RowHeight = 31.5

# --- Selection / active cell ----------------------------------------------
$ws.Range("D12").Select()
